# Scheduled runner update: refresh market-derived price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the ALC, ARM,
# CUL, LTW and WVR leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 222.16667
$ws.Range("I8").Value = 70.2
$ws.Range("J8").Value = 982
$ws.Range("K8").Value = 210.6
$ws.Range("L8").Value = 2946
$ws.Range("M8").Value = -71.60000000000002
$ws.Range("N8").Value = -3224
$ws.Range("H64").Value = 3158.9395
$ws.Range("J64").Value = 3184.5454
$ws.Range("L64").Value = 3184.5454
$ws.Range("N64").Value = -3680.5454
$ws.Range("H67").Value = 3158.9395
$ws.Range("J67").Value = 3184.5454
$ws.Range("L67").Value = 3184.5454
$ws.Range("N67").Value = -4900.5454
$ws.Range("H68").Value = 29095.5
$ws.Range("J68").Value = 29095.5
$ws.Range("L68").Value = 29095.5
$ws.Range("N68").Value = -30593.5
$ws.Range("H69").Value = 5492.75
$ws.Range("I69").Value = 2106.5
$ws.Range("J69").Value = 6170
$ws.Range("K69").Value = 6319.5
$ws.Range("L69").Value = 18510
$ws.Range("M69").Value = -5445.5
$ws.Range("N69").Value = -20258
$ws.Range("H70").Value = 53334390
$ws.Range("I70").Value = 93333890
$ws.Range("J70").Value = 1733
$ws.Range("K70").Value = 280001670
$ws.Range("L70").Value = 5199
$ws.Range("M70").Value = -280001400
$ws.Range("N70").Value = -5739
$ws.Range("H71").Value = 29095.5
$ws.Range("J71").Value = 29095.5
$ws.Range("L71").Value = 87286.5
$ws.Range("N71").Value = -94774.5
$ws.Range("H72").Value = 5492.75
$ws.Range("I72").Value = 2106.5
$ws.Range("J72").Value = 6170
$ws.Range("K72").Value = 18958.5
$ws.Range("L72").Value = 55530
$ws.Range("M72").Value = -14590.5
$ws.Range("N72").Value = -64266
$ws.Range("H73").Value = 53334390
$ws.Range("I73").Value = 93333890
$ws.Range("J73").Value = 1733
$ws.Range("K73").Value = 280001670
$ws.Range("L73").Value = 5199
$ws.Range("M73").Value = -280000734
$ws.Range("N73").Value = -7071
$ws.Range("H74").Value = 6050.8887
$ws.Range("I74").Value = 6642.857
$ws.Range("K74").Value = 6642.857
$ws.Range("M74").Value = -5706.857
$ws.Range("H75").Value = 27434.666
$ws.Range("J75").Value = 27434.666
$ws.Range("L75").Value = 27434.666
$ws.Range("N75").Value = -29306.666
$ws.Range("H76").Value = 3597.25
$ws.Range("I76").Value = 3597.25
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3597.25
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3282.25
$ws.Range("N76").ClearContents()
$ws.Range("H77").Value = 6050.8887
$ws.Range("I77").Value = 6642.857
$ws.Range("K77").Value = 33214.285
$ws.Range("M77").Value = -28534.285
$ws.Range("H78").Value = 27434.666
$ws.Range("J78").Value = 27434.666
$ws.Range("L78").Value = 82303.99800000001
$ws.Range("N78").Value = -91663.99800000001
$ws.Range("H79").Value = 3597.25
$ws.Range("I79").Value = 3597.25
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3597.25
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2505.25
$ws.Range("N79").ClearContents()
$ws.Range("H80").Value = 593.0769
$ws.Range("I80").Value = 296.66666
$ws.Range("J80").Value = 750
$ws.Range("K80").Value = 889.9999799999999
$ws.Range("L80").Value = 2250
$ws.Range("M80").Value = 108.0000200000001
$ws.Range("N80").Value = -4246
$ws.Range("H82").Value = 833.3333
$ws.Range("I82").Value = 833.3333
$ws.Range("K82").Value = 2499.9999
$ws.Range("M82").Value = -2093.9999
$ws.Range("H83").Value = 593.0769
$ws.Range("I83").Value = 296.66666
$ws.Range("J83").Value = 750
$ws.Range("K83").Value = 2669.99994
$ws.Range("L83").Value = 6750
$ws.Range("M83").Value = 2322.00006
$ws.Range("N83").Value = -16734
$ws.Range("H85").Value = 833.3333
$ws.Range("I85").Value = 833.3333
$ws.Range("K85").Value = 2499.9999
$ws.Range("M85").Value = -1095.9999
$ws.Range("H86").Value = 2322.2222
$ws.Range("J86").Value = 2833.3333
$ws.Range("L86").Value = 2833.3333
$ws.Range("N86").Value = -5079.3333
$ws.Range("H87").Value = 35132
$ws.Range("J87").Value = 35132
$ws.Range("L87").Value = 35132
$ws.Range("N87").Value = -37628
$ws.Range("H88").Value = 348
$ws.Range("I88").Value = 348
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 348
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 58
$ws.Range("N88").ClearContents()
$ws.Range("H89").Value = 2322.2222
$ws.Range("J89").Value = 2833.3333
$ws.Range("L89").Value = 14166.6665
$ws.Range("N89").Value = -25398.6665
$ws.Range("H90").Value = 35132
$ws.Range("J90").Value = 35132
$ws.Range("L90").Value = 105396
$ws.Range("N90").Value = -117876
$ws.Range("H91").Value = 348
$ws.Range("I91").Value = 348
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 348
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 1056
$ws.Range("N91").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 19571.428
$ws.Range("I134").Value = 9000
$ws.Range("J134").Value = 27500
$ws.Range("K134").Value = 9000
$ws.Range("L134").Value = 27500
$ws.Range("M134").Value = -3930
$ws.Range("N134").Value = -37640

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 872.8
$ws.Range("I11").Value = 134
$ws.Range("J11").Value = 2350.4
$ws.Range("K11").Value = 402
$ws.Range("L11").Value = 7051.200000000001
$ws.Range("M11").Value = -262
$ws.Range("N11").Value = -7331.200000000001
$ws.Range("H131").Value = 831.4792
$ws.Range("J131").Value = 884.7907
$ws.Range("L131").Value = 2654.3721
$ws.Range("N131").Value = -12734.3721

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1599.5454
$ws.Range("I68").Value = 1555
$ws.Range("K68").Value = 1555
$ws.Range("M68").Value = -806
$ws.Range("H71").Value = 1599.5454
$ws.Range("I71").Value = 1555
$ws.Range("K71").Value = 7775
$ws.Range("M71").Value = -4031

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 3250
$ws.Range("J62").Value = 4500
$ws.Range("K62").Value = 3250
$ws.Range("L62").Value = 4500
$ws.Range("M62").Value = -2626
$ws.Range("N62").Value = -5748
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 3250
$ws.Range("J65").Value = 4500
$ws.Range("K65").Value = 16250
$ws.Range("L65").Value = 22500
$ws.Range("M65").Value = -13130
$ws.Range("N65").Value = -28740
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H81").Value = 13722.111
$ws.Range("I81").Value = 50247.5
$ws.Range("J81").Value = 3286.2856
$ws.Range("K81").Value = 100495
$ws.Range("L81").Value = 6572.5712
$ws.Range("M81").Value = -99434
$ws.Range("N81").Value = -8694.5712
$ws.Range("H84").Value = 13722.111
$ws.Range("I84").Value = 50247.5
$ws.Range("J84").Value = 3286.2856
$ws.Range("K84").Value = 502475
$ws.Range("L84").Value = 32862.856
$ws.Range("M84").Value = -497171
$ws.Range("N84").Value = -43470.856
